$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with rich-text runs) ---
$ws.Range("A8").Value = "Volume 32   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/28/2025  Through  8/3/2025"

# --- Cells changing from numeric to text placeholder ("0" / "***.*") ---
$ws.Range("C14").Copy($ws.Range("D14"))
$ws.Range("M14").Copy($ws.Range("E14"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("M14").Copy($ws.Range("E22"))

# --- Cells changing from text placeholder to numeric ---
$ws.Range("D16").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("D16").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 1
$ws.Range("E16").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 100
$ws.Range("D16").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("E16").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("D16").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("E16").Copy($ws.Range("H31"))
$ws.Range("H31").Value = -100

# --- Simple numeric value updates ---
$ws.Range("M15").Value = 75
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 11.111111111111
$ws.Range("I16").Value = 62
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = -7.462686567164
$ws.Range("L16").Value = 31.914893617021
$ws.Range("M16").Value = -37.373737373737
$ws.Range("N16").Value = -85.514018691588
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 38.461538461538
$ws.Range("I17").Value = 132
$ws.Range("J17").Value = 124
$ws.Range("K17").Value = 6.451612903225
$ws.Range("L17").Value = 3.125
$ws.Range("M17").Value = 48.314606741573
$ws.Range("N17").Value = -30.89005235602
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 79
$ws.Range("J18").Value = 66
$ws.Range("K18").Value = 19.696969696969
$ws.Range("L18").Value = 14.492753623188
$ws.Range("M18").Value = -63.926940639269
$ws.Range("N18").Value = -92.004048582995
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -41.666666666666
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 40
$ws.Range("H19").Value = -20
$ws.Range("I19").Value = 256
$ws.Range("J19").Value = 323
$ws.Range("K19").Value = -20.743034055727
$ws.Range("L19").Value = -22.424242424242
$ws.Range("M19").Value = 11.304347826087
$ws.Range("N19").Value = -32.453825857519
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 17
$ws.Range("H20").Value = 13.333333333333
$ws.Range("I20").Value = 133
$ws.Range("J20").Value = 107
$ws.Range("K20").Value = 24.29906542056
$ws.Range("L20").Value = 47.777777777777
$ws.Range("M20").Value = 54.651162790697
$ws.Range("N20").Value = -87.876025524156
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -1.123595505617
$ws.Range("I21").Value = 669
$ws.Range("J21").Value = 703
$ws.Range("K21").Value = -4.836415362731
$ws.Range("L21").Value = -1.035502958579
$ws.Range("M21").Value = -7.977991746905
$ws.Range("N21").Value = -78.440219142765
$ws.Range("M22").Value = -25
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 65
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = -9.722222222222
$ws.Range("I24").Value = 586
$ws.Range("J24").Value = 579
$ws.Range("K24").Value = 1.208981001727
$ws.Range("L24").Value = -11.346444780635
$ws.Range("M24").Value = 12.692307692307
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -60
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -57.894736842105
$ws.Range("I25").Value = 92
$ws.Range("J25").Value = 142
$ws.Range("K25").Value = -35.211267605633
$ws.Range("L25").Value = -52.820512820512
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 12
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = -45.454545454545
$ws.Range("I26").Value = 222
$ws.Range("J26").Value = 244
$ws.Range("K26").Value = -9.016393442622
$ws.Range("L26").Value = 6.730769230769
$ws.Range("M26").Value = -15.589353612167
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 11
$ws.Range("K27").Value = -47.619047619047
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 2
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -42.857142857142
$ws.Range("I28").Value = 34
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = -15
$ws.Range("L28").Value = -15
$ws.Range("J31").Value = 11
$ws.Range("K31").Value = -27.272727272727
